# Perfect Doctor workbook update:
#  - Insert a new worksheet "Sheet1" between "Sheet3" and "Sheet2" with a
#    fresh copy of the person/Y1/Y0 table plus the step-by-step exercise
#    prompts.
#  - Tweak the saved selections on "Sheet3" and (the now renamed) "Sheet2"
#    sheet views.

$wb = $excel.ActiveWorkbook

$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet2 = $wb.Worksheets.Item("Sheet2")

# --- Insert the new worksheet right before "Sheet2" -------------------
$new = $wb.Worksheets.Add($sheet2)
$new.Name = "Sheet1"

# Re-fetch the other sheet references: the underlying collection shifted
# when the new sheet was inserted, so the old handles can go stale.
$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet2 = $wb.Worksheets.Item("Sheet2")
$new = $wb.Worksheets.Item("Sheet1")

# Reuse the existing table formatting (fills / number formats / borders)
# from the "Sheet3" tab so the new sheet's styles line up with the rest
# of the workbook instead of creating brand-new style entries.
$sheet3.Range("A1:F1").Copy()
$new.Range("A1:F1").PasteSpecial(-4122)

$sheet3.Range("A2:D13").Copy()
$new.Range("A2:D13").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Header row ---------------------------------------------------------
$new.Range("A1").Value = "Person"
$new.Range("B1").Value = "Y0"
$new.Range("C1").Value = "Y1"
$new.Range("D1").Value = "TE"
$new.Range("E1").Value = "D"
$new.Range("F1").Value = "Y"

# --- Person / Y0 / Y1 table ---------------------------------------------
$people = @("Adam","Billy","Cynthia","Daniel","Elaine","Francis","Gia","Hank","Ida","Jane","Kelly","Leanna")
$y0 = @(48,27,33,33,50,26,36,35,46,33,38,47)
$y1 = @(25,36,33,31,55,38,29,23,34,46,27,24)

for ($i = 0; $i -lt $people.Length; $i++) {
    $row = 2 + $i
    $new.Range("A$row").Value = $people[$i]
    $new.Range("B$row").Value = $y0[$i]
    $new.Range("C$row").Value = $y1[$i]
}

# --- Exercise step prompts ------------------------------------------------
$new.Range("A15").Value = "Step 1: Calculate TE"
$new.Range("A16").Value = "Step 2: Use perfect doctor to assign D = 1 if TE>0 otherwise D=0"
$new.Range("A17").Value = "Step 3: Use switching equation to get Y"
$new.Range("A18").Value = "Step 4: Calculate SDO"
$new.Range("A19").Value = "Step 5: Calculate selection bias, pi, ATE, ATT, ATU"
$new.Range("A20").Value = "Step 6: Show the decomposition of the SDO equals the sum of ATE, selection bias and heterogenous treatment effects bias"
$new.Range("A22").Value = "How much of the SDO is due to ""causal effect"" and how much is due to ""selection bias"""
$new.Range("A23").Value = "What is selection bias now that you have gone through this exercise? Put into words that you could tell your parent and they would understand"

# --- Sheet view tweaks ----------------------------------------------------
$sheet3.Activate()
$sheet3.Range("A1:F13").Select() | Out-Null

$sheet2.Activate()
$sheet2.Range("B16").Select() | Out-Null

$new.Activate()
$new.Range("A24").Select() | Out-Null
$excel.ActiveWindow.Zoom = 280
